$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D value (only rows whose Price column changed).
# Values that look like plain decimal numbers get a leading apostrophe so
# Excel keeps them as text (matching the source inlineStr cells) instead of
# silently converting them to numeric cells. Values with two or more dots
# (e.g. "61.676.12") already fail numeric parsing, so no apostrophe is
# needed for those - adding one would needlessly flag them with a
# quote-prefix cell style that the original file doesn't have.
$dUpdates = @{
    2  = "61.676.12"
    3  = "3.045.66"
    5  = "'580.61"
    6  = "'130.17"
    8  = "3.039.64"
    11 = "'5.22"
    13 = "'0.0000231"
    14 = "'33.38"
    16 = "3.542.47"
    17 = "61.619.18"
    18 = "3.049.70"
    20 = "'447.41"
    25 = "'12.81"
    26 = "'0.999"
    30 = "'7.42"
    32 = "'25.89"
    33 = "'0.0967"
    35 = "'0.968"
    37 = "'50.29"
    39 = "'0.0372"
    40 = "'7.91"
    42 = "'2.53"
    43 = "'377.43"
    44 = "2.686.01"
    45 = "'0.999"
    46 = "'122.90"
    51 = "'23.76"
}

# Map of row -> new E value (Volume(1h) column), all rows 2..51 change
$eUpdates = @{
    2  = "  -1.43%  "
    3  = "  -4.40%  "
    4  = "  -0.13%  "
    5  = "  -1.21%  "
    6  = "  -4.26%  "
    7  = "  +0.00%  "
    8  = "  -4.47%  "
    9  = "  -1.18%  "
    10 = "  -2.78%  "
    11 = "  -0.69%  "
    12 = "  -3.56%  "
    13 = "  -1.55%  "
    14 = "  +0.12%  "
    15 = "  +0.93%  "
    16 = "  -4.63%  "
    17 = "  -1.56%  "
    18 = "  -4.35%  "
    19 = "  -2.40%  "
    20 = "  -1.95%  "
    21 = "  -3.32%  "
    22 = "  -4.74%  "
    23 = "  -4.00%  "
    24 = "  -3.37%  "
    25 = "  -3.10%  "
    26 = "  -0.12%  "
    27 = "  -0.19%  "
    28 = "  -4.97%  "
    29 = "  -1.00%  "
    30 = "  -4.56%  "
    31 = "  -5.28%  "
    32 = "  -5.29%  "
    33 = "  -6.41%  "
    34 = "  -2.16%  "
    35 = "  -6.88%  "
    36 = "  -3.38%  "
    37 = "  -1.57%  "
    38 = "  -0.01%  "
    39 = "  -3.17%  "
    40 = "  -1.09%  "
    41 = "  -2.00%  "
    42 = "  -7.04%  "
    43 = "  -3.41%  "
    44 = "  -5.24%  "
    45 = "  +0.02%  "
    46 = "  -1.65%  "
    47 = "  -4.36%  "
    48 = "  -5.87%  "
    49 = "  -5.88%  "
    50 = "  -2.78%  "
    51 = "  -6.69%  "
}

foreach ($row in $dUpdates.Keys) {
    $ws.Range("D$row").Value = $dUpdates[$row]
}

foreach ($row in $eUpdates.Keys) {
    $ws.Range("E$row").Value = $eUpdates[$row]
}
